$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.65903648854455
$ws.Range("C2").Value = 7.932794846617679
$ws.Range("E2").Value = 11.882420359395
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 48.80867767209238
$ws.Range("H2").Value = 19.25876126438727
$ws.Range("K2").Value = 10.98704373276778
$ws.Range("L2").Value = 9.874170581194372
$ws.Range("B3").Value = 14.48795770986589
$ws.Range("C3").Value = 7.914179842343547
$ws.Range("E3").Value = 11.8855654352398
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 48.7675845688211
$ws.Range("H3").Value = 19.29387742554934
$ws.Range("K3").Value = 10.8777593076427
$ws.Range("L3").Value = 9.866265104004345
$ws.Range("B4").Value = 14.38565901498351
$ws.Range("C4").Value = 7.902440111922375
$ws.Range("E4").Value = 11.88953797794296
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 48.75407083715752
$ws.Range("H4").Value = 19.31843298937285
$ws.Range("K4").Value = 10.81276876696415
$ws.Range("L4").Value = 9.863190286126908
$ws.Range("B5").Value = 14.34471127972582
$ws.Range("C5").Value = 7.897578465919604
$ws.Range("E5").Value = 11.8916705669372
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 48.75150749853966
$ws.Range("H5").Value = 19.32919107445932
$ws.Range("K5").Value = 10.78684456307202
$ws.Range("L5").Value = 9.862386091414113
$ws.Range("B6").Value = 14.33795804668013
$ws.Range("C6").Value = 7.896766516167235
$ws.Range("E6").Value = 11.89205571888647
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 48.75125951745874
$ws.Range("H6").Value = 19.33102279214914
$ws.Range("K6").Value = 10.782574523207
$ws.Range("L6").Value = 9.862279699910642
$ws.Range("B7").Value = 14.38510371943111
$ws.Range("C7").Value = 7.90237485983227
$ws.Range("E7").Value = 11.88956465824659
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 48.75402435405332
$ws.Range("H7").Value = 19.31857503585243
$ws.Range("K7").Value = 10.81241683895256
$ws.Range("L7").Value = 9.863177621523127
$ws.Range("B8").Value = 14.59950918234689
$ws.Range("C8").Value = 7.926440949688591
$ws.Range("E8").Value = 11.88308137307036
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 48.79207571122453
$ws.Range("H8").Value = 19.27024707454516
$ws.Range("K8").Value = 10.94894229118964
$ws.Range("L8").Value = 9.87107643393751
$ws.Range("B9").Value = 15.03948759085106
$ws.Range("C9").Value = 7.971160416197695
$ws.Range("E9").Value = 11.88654162442651
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 48.95970304884612
$ws.Range("H9").Value = 19.19929316578101
$ws.Range("K9").Value = 11.23208183483356
$ws.Range("L9").Value = 9.900610590011029
$ws.Range("B10").Value = 15.37158211422946
$ws.Range("C10").Value = 8.002490906037536
$ws.Range("E10").Value = 11.89890234799
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 49.1394718453048
$ws.Range("H10").Value = 19.16176192320129
$ws.Range("K10").Value = 11.44766478171284
$ws.Range("L10").Value = 9.93076033971494
$ws.Range("B11").Value = 15.52392440489634
$ws.Range("C11").Value = 8.016409447351554
$ws.Range("E11").Value = 11.90664468422956
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 49.23347890469099
$ws.Range("H11").Value = 19.14787306960611
$ws.Range("K11").Value = 11.54698144034388
$ws.Range("L11").Value = 9.946282431697394
$ws.Range("B12").Value = 15.58174075669224
$ws.Range("C12").Value = 8.021631693435269
$ws.Range("E12").Value = 11.90987979621856
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 49.27082563616187
$ws.Range("H12").Value = 19.14307267158235
$ws.Range("K12").Value = 11.58473577388228
$ws.Range("L12").Value = 9.952417171078093
$ws.Range("B13").Value = 15.56928416290105
$ws.Range("C13").Value = 8.020509148324905
$ws.Range("E13").Value = 11.90916959575025
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 49.26270478216084
$ws.Range("H13").Value = 19.144086092094
$ws.Range("K13").Value = 11.57659878135533
$ws.Range("L13").Value = 9.951084571065845
$ws.Range("B14").Value = 15.52867874571719
$ws.Range("C14").Value = 8.016840057010826
$ws.Range("E14").Value = 11.90690477049013
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 49.23651647157409
$ws.Range("H14").Value = 19.14746893276995
$ws.Range("K14").Value = 11.55008479486485
$ws.Range("L14").Value = 9.946782010385562
$ws.Range("B15").Value = 15.50382175512464
$ws.Range("C15").Value = 8.014586317398088
$ws.Range("E15").Value = 11.90555694436455
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 49.22070270966719
$ws.Range("H15").Value = 19.14960082802491
$ws.Range("K15").Value = 11.53386212433495
$ws.Range("L15").Value = 9.94417992636842
$ws.Range("B16").Value = 15.36164703086848
$ws.Range("C16").Value = 8.001574544733277
$ws.Range("E16").Value = 11.89843889143412
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 49.1335735256807
$ws.Range("H16").Value = 19.16273376616028
$ws.Range("K16").Value = 11.44119638233466
$ws.Range("L16").Value = 9.929782045238071
$ws.Range("B17").Value = 15.27471348721416
$ws.Range("C17").Value = 7.993506452952922
$ws.Range("E17").Value = 11.89461401518874
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 49.08324827908579
$ws.Range("H17").Value = 19.17160679884344
$ws.Range("K17").Value = 11.38464405243394
$ws.Range("L17").Value = 9.921410157354606
$ws.Range("B18").Value = 15.22483448823579
$ws.Range("C18").Value = 7.988834499377837
$ws.Range("E18").Value = 11.89261363724034
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 49.05545437853749
$ws.Range("H18").Value = 19.17700999394017
$ws.Range("K18").Value = 11.35223604445276
$ws.Range("L18").Value = 9.916765128627597
$ws.Range("B19").Value = 15.20796911185764
$ws.Range("C19").Value = 7.987247275179911
$ws.Range("E19").Value = 11.89197066584228
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 49.04624194627025
$ws.Range("H19").Value = 19.17889085447291
$ws.Range("K19").Value = 11.3412847956911
$ws.Range("L19").Value = 9.915221729962715
$ws.Range("B20").Value = 15.2839553958318
$ws.Range("C20").Value = 7.994368567820062
$ws.Range("E20").Value = 11.89500053326523
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 49.08848634196636
$ws.Range("H20").Value = 19.17063122709838
$ws.Range("K20").Value = 11.39065203446699
$ws.Range("L20").Value = 9.922283758491465
$ws.Range("B21").Value = 15.54060251219893
$ws.Range("C21").Value = 8.017919074862077
$ws.Range("E21").Value = 11.90756178730502
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 49.24416125856231
$ws.Range("H21").Value = 19.14646284476881
$ws.Range("K21").Value = 11.557868923024
$ws.Range("L21").Value = 9.948038830477101
$ws.Range("B22").Value = 15.70905526900684
$ws.Range("C22").Value = 8.033028379076335
$ws.Range("E22").Value = 11.91753806148398
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 49.35608779313074
$ws.Range("H22").Value = 19.13334330459191
$ws.Range("K22").Value = 11.66798597029674
$ws.Range("L22").Value = 9.966366886101925
$ws.Range("B23").Value = 15.61910104825184
$ws.Range("C23").Value = 8.024990220596665
$ws.Range("E23").Value = 11.91205242881137
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 49.29542266854518
$ws.Range("H23").Value = 19.14010026721186
$ws.Range("K23").Value = 11.6091495617881
$ws.Range("L23").Value = 9.956449050299444
$ws.Range("B24").Value = 15.27977681243945
$ws.Range("C24").Value = 7.993978910023462
$ws.Range("E24").Value = 11.89482516973841
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 49.08611466571298
$ws.Range("H24").Value = 19.17107134248394
$ws.Range("K24").Value = 11.38793549710917
$ws.Range("L24").Value = 9.921888279643614
$ws.Range("B25").Value = 14.91868563426553
$ws.Range("C25").Value = 7.959329206594617
$ws.Range("E25").Value = 11.88387749562232
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 48.90439782184707
$ws.Range("H25").Value = 19.21592975133762
$ws.Range("K25").Value = 11.15402247309554
$ws.Range("L25").Value = 9.891127671410267
